$wb = $excel.ActiveWorkbook

# Delete Sheet2 (this removes its associated drawing/images too)
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete()
$excel.DisplayAlerts = $true

# Rename Sheet1
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "PKES relay replay attack"

# Update the view: topLeftCell and selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3  # C
$ws.Range("H7").Select()
